$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert rows 4, 5, 6 (copy row 3 pattern: A/E/F/J/K, styles 1/1/1/4/4) ---
$ws.Rows.Item(3).Copy()
$ws.Rows.Item(4).Insert()

$ws.Rows.Item(3).Copy()
$ws.Rows.Item(5).Insert()

$ws.Rows.Item(3).Copy()
$ws.Rows.Item(6).Insert()

# --- Row 4: Thriveni.docx ---
$edu4 = @"
[
    {
        "college": "KITS College, Kodad, Telangana",
        "degree": "B.Tech",
        "branch": "",
        "yop": "",
        "cgpa/%": "",
        "isFinal": 0
    }
]
"@
$work4 = @"
[
    {
        "company": "UI Sottech Pvt Ltd",
        "designation": "",
        "fromMonth": "07",
        "fromYear": "2022",
        "toMonth": "07",
        "toYear": "2024",
        "IsLatest": 1
    }
]
"@
$ws.Range("A4").Value = "Thriveni.docx"
$ws.Range("E4").Value = "None"
$ws.Range("F4").Value = "Unknown"
$ws.Range("J4").Value = $edu4
$ws.Range("K4").Value = $work4
$ws.Rows.Item(4).RowHeight = 127.85

# --- Row 5: Vipin_yadav_cv.pdf ---
$edu5 = @"
[
    {
        "college": "Thakur tej Bahadur institute of technology karampur saidpur Ghazipur",
        "degree": "Diploma",
        "branch": "Electrical Engineering",
        "yop": "2020",
        "cgpa/%": "",
        "isFinal": 1
    },
    {
        "college": "Up Board",
        "degree": "12th",
        "branch": "",
        "yop": "2017",
        "cgpa/%": "79",
        "isFinal": 0
    },
    {
        "college": "Up board",
        "degree": "10th",
        "branch": "",
        "yop": "2015",
        "cgpa/%": "83",
        "isFinal": 0
    }
]
"@
$work5 = @"
[
    {
        "company": "Samvdhana Motherson internation ltd.",
        "designation": "Assistant Engineer",
        "fromMonth": "06",
        "fromYear": "2022",
        "toMonth": "06",
        "toYear": "25",
        "IsLatest": 1
    }
]
"@
$ws.Range("A5").Value = "Vipin_yadav_cv.pdf"
$ws.Range("E5").Value = "Male"
$ws.Range("F5").Value = "Single"
$ws.Range("G5").Value = "20/07/2000"
$ws.Range("J5").Value = $edu5
$ws.Range("K5").Value = $work5
$ws.Rows.Item(5).RowHeight = 299.5

# --- Row 6: ArreeshRajan.docx ---
$edu6 = @"
[
    {
        "college": "",
        "degree": "Executive MBA",
        "branch": "",
        "yop": "",
        "cgpa/%": "",
        "isFinal": 0
    },
    {
        "college": "",
        "degree": "B.Tech",
        "branch": "Information Technology",
        "yop": "",
        "cgpa/%": "",
        "isFinal": 0
    }
]
"@
$ws.Range("A6").Value = "ArreeshRajan.docx"
$ws.Range("E6").Value = "None"
$ws.Range("F6").Value = "Unknown"
$ws.Range("J6").Value = $edu6
$ws.Rows.Item(6).RowHeight = 207.95

# --- Column K width ---
$ws.Columns.Item(11).ColumnWidth = 51.96

# --- Sheet view: freeze top row, scroll to H1/K6 selection ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("H1").Select()
$ws.Range("K6").Select()
